# Actualización automática de grupos experimentales
# Swap "Sin SmartScore" <-> "Con SmartScore" labels in column B for the
# affected participant rows, fill in the previously-missing label for
# row 27, and fix the SmartScore numeric cells in row 27 that were
# incorrectly stored as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Grupo_Experimental (column B) label needs to be swapped
$rowsSinToCon = @(3, 7, 15, 18, 19, 26)
$rowsConToSin = @(9, 12, 16, 21, 22, 24)

foreach ($r in $rowsSinToCon) {
    $ws.Range("B$r").Value = "Con SmartScore"
}

foreach ($r in $rowsConToSin) {
    $ws.Range("B$r").Value = "Sin SmartScore"
}

# Row 27 previously had an empty Grupo_Experimental cell
$ws.Range("B27").Value = "Sin SmartScore"

# Row 27 SmartScore values were stored as text; convert them to real numbers
$ws.Range("I27").Value = 0.578
$ws.Range("L27").Value = 0.566
$ws.Range("O27").Value = 0.455
$ws.Range("R27").Value = 0.712
$ws.Range("U27").Value = 0.625
$ws.Range("X27").Value = 0.567
$ws.Range("AA27").Value = 0.657
$ws.Range("AD27").Value = 0.656
$ws.Range("AG27").Value = 0.644
